$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# List2 ("List2" / sheet2.xml): add a new "Obtížnost" (difficulty) column
# between the existing "Proher" (C) and "V kolikátém kole" columns, and fill
# in a bunch of previously-empty B/C/D/E/F values.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("List2")

# New header for the inserted column, and move the "V kolikátém kole" header
# from E1 to F1 (write the new location first so the shared string stays
# referenced when we clear the old cell).
$ws2.Range("D1").Value = "Obtížnost"
$ws2.Range("F1").Value = "V kolikátém kole"
$ws2.Range("E1").ClearContents() | Out-Null

# Row 2 - Prvni kroky
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 1

# Row 3 - Trziste otevrena
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 37

# Row 4 - Neznamy kolonizator
$ws2.Range("C4").Value = 0
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 28

# Row 5 - Poustni mesto
$ws2.Range("B5").Value = 0
$ws2.Range("D5").Value = 2

# Row 6 - Mnisi a kupci
$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 2
$ws2.Range("F6").Value = 31

# Row 7 - Rybi ostrov
$ws2.Range("B7").Value = 1
$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 1
$ws2.Range("E7").Value = 4
$ws2.Range("F7").Value = 28

# Row 8 - Tri ostruvky
$ws2.Range("B8").Value = 1
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 1
$ws2.Range("E8").Value = 5
$ws2.Range("F8").Value = 35

# Row 9 - Posvatna hora
$ws2.Range("B9").Value = 0
$ws2.Range("C9").Value = 1

# Row 10 - Boj o moc
$ws2.Range("B10").Value = 0
$ws2.Range("C10").Value = 3
$ws2.Range("D10").Value = 3

# Row 11 - Okruzni cesta
$ws2.Range("B11").Value = 1
$ws2.Range("C11").Value = 3
$ws2.Range("D11").Value = 3
$ws2.Range("F11").Value = 67

# Row 12 - Zlaty dul
$ws2.Range("B12").Value = 1
$ws2.Range("C12").Value = 0
$ws2.Range("F12").Value = 71

# Row 13 - Nedostanek kamene
$ws2.Range("B13").Value = 1
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 2
$ws2.Range("F13").Value = 69

# Row 14 - Mec a stit
$ws2.Range("B14").Value = 1
$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 1
$ws2.Range("F14").Value = 59

# Row 15 - Pan cest
$ws2.Range("B15").Value = 1
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 1
$ws2.Range("F15").Value = 42

# Row 16 - Temne oci
$ws2.Range("B16").Value = 1
$ws2.Range("C16").Value = 0
$ws2.Range("D16").Value = 2
$ws2.Range("F16").Value = 48

# ---------------------------------------------------------------------------
# List3 ("List3" / sheet3.xml): fill in a little running-total table with
# differences computed via a (shared) formula.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("List3")

$ws3.Range("C2").Value = 691
$ws3.Range("C3").Value = 802
$ws3.Range("C4").Value = 892
$ws3.Range("C5").Value = 1000
$ws3.Range("C6").Value = 1178
$ws3.Range("C7").Value = 1334
$ws3.Range("C8").Value = 1449

$ws3.Range("D3").Formula = "=C3-C2"
$ws3.Range("D4:D8").Formula = "=C4-C3"

# First value is emphasised (bold), matching the existing bold style already
# used elsewhere in the workbook.
$ws3.Range("C2").Font.Bold = $true

# ---------------------------------------------------------------------------
# Selections / active sheet: List3 keeps a lingering selection at A8, but the
# workbook's active tab ends up on List2 with D13 selected (selecting List2
# last is what flips the active tab / tabSelected flags correctly).
# ---------------------------------------------------------------------------
$ws3.Range("A8").Select() | Out-Null
$ws2.Range("D13").Select() | Out-Null
